# Update the "取得日時" (retrieved at) timestamp for the newly (re-)appended
# rows on the "ランサーズ" sheet from 2025-11-13 06:28:54 to 2025-11-13 06:35:39.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-13 06:35:39"

for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
